$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.114.67"
$ws.Range("E2").Value = "  -2.65%  "
$ws.Range("D3").Value = "3.452.57"
$ws.Range("E3").Value = "  -2.52%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "575.34"
$ws.Range("E5").Value = "  +1.11%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "175.38"
$ws.Range("E6").Value = "  -7.38%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.628"
$ws.Range("E7").Value = "  +1.39%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.629"
$ws.Range("E9").Value = "  -1.12%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.159"
$ws.Range("E10").Value = "  +4.80%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "55.59"
$ws.Range("E11").Value = "  +1.19%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000275"
$ws.Range("E12").Value = "  +1.46%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "9.18"
$ws.Range("E13").Value = "  -2.65%  "
$ws.Range("D14").Value = "4.009.16"
$ws.Range("E14").Value = "  -2.49%  "
$ws.Range("E15").Value = "  -0.57%  "
$ws.Range("D16").Value = "3.457.71"
$ws.Range("E16").Value = "  -2.40%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "18.22"
$ws.Range("E17").Value = "  -0.30%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "11.97"
$ws.Range("E18").Value = "  -0.65%  "
$ws.Range("D19").Value = "65.201.32"
$ws.Range("E19").Value = "  -2.54%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.998"
$ws.Range("E20").Value = "  -0.14%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "408.35"
$ws.Range("E21").Value = "  -5.29%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.23"
$ws.Range("E22").Value = "  +0.07%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.37"
$ws.Range("E23").Value = "  +5.89%  "
$ws.Range("B24").Value = "InternetComputer(DFINITY)"
$ws.Range("C24").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.49"
$ws.Range("E24").Value = "  +9.70%  "
$ws.Range("B25").Value = "Litecoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "84.16"
$ws.Range("E25").Value = "  -1.21%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.90"
$ws.Range("E26").Value = "  -2.05%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.83"
$ws.Range("E27").Value = "  -2.97%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.08"
$ws.Range("E28").Value = "  -1.69%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "30.02"
$ws.Range("E29").Value = "  -1.28%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.65"
$ws.Range("E30").Value = "  +0.24%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "11.61"
$ws.Range("E31").Value = "  -1.44%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "588.17"
$ws.Range("E32").Value = "  -8.53%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.109"
$ws.Range("E33").Value = "  -3.02%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "59.95"
$ws.Range("E34").Value = "  +0.14%  "
$ws.Range("B35").Value = "Dai"
$ws.Range("C35").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.00"
$ws.Range("E35").Value = "  +0.25%  "
$ws.Range("B36").Value = "Kaspa"
$ws.Range("C36").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.154"
$ws.Range("E36").Value = "  +1.88%  "
$ws.Range("D37").Value = "0.0₃0781"
$ws.Range("E37").Value = "  -4.41%  "
$ws.Range("B38").Value = "Stacks"
$ws.Range("C38").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.55"
$ws.Range("E38").Value = "  +5.99%  "
$ws.Range("B39").Value = "InjectiveProtocol"
$ws.Range("C39").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "36.43"
$ws.Range("E39").Value = "  -5.98%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.379"
$ws.Range("E40").Value = "  -3.25%  "
$ws.Range("D41").Value = "3.207.23"
$ws.Range("E41").Value = "  +4.84%  "
$ws.Range("E42").Value = "  -0.08%  "
$ws.Range("E43").Value = "  +2.45%  "
$ws.Range("B44").Value = "Fetch.AI"
$ws.Range("C44").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.52"
$ws.Range("E44").Value = "  -5.99%  "
$ws.Range("B45").Value = "ApeXProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.27"
$ws.Range("E45").Value = "  -2.43%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0414"
$ws.Range("E46").Value = "  -1.67%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.131"
$ws.Range("E47").Value = "  -0.76%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.63"
$ws.Range("E48").Value = "  -6.09%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.55"
$ws.Range("E49").Value = "  -1.16%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "137.74"
$ws.Range("E50").Value = "  -2.58%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.34"
$ws.Range("E51").Value = "  -2.73%  "
